$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.097.21'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +4.92%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.789.41'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +7.75%  '

# Row 4
$ws.Range("E4").Value = '  -0.13%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '428.21'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +10.46%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.00'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +16.45%  '

# Row 7
$ws.Range("E7").Value = '  +5.91%  '

# Row 8
$ws.Range("E8").Value = '  -0.03%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.740'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +10.43%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.156'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.13%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000325'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.32%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '43.53'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +13.75%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.69'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +18.45%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.379.59'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +7.68%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.00'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +17.38%  '

# Row 16
$ws.Range("E16").Value = '  +1.53%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.781.69'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +7.28%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '20.17'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +9.07%  '

# Row 19
$ws.Range("E19").Value = '  +13.21%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '66.268.89'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.14%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '411.47'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.99%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '15.21'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +10.18%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.31'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +16.00%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.65'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +6.14%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '37.13'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +11.65%  '

# Row 26
$ws.Range("B26").Value = 'PancakeSwap'
$ws.Range("C26").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.29'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +12.11%  '

# Row 27
$ws.Range("B27").Value = 'RenderToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.74'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +47.06%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.84'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +15.07%  '

# Row 29
$ws.Range("E29").Value = '  -0.79%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '13.97'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +20.27%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '707.26'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +7.31%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.130'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +18.86%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.79'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.62%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '40.29'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +11.72%  '

# Row 35
$ws.Range("B35").Value = 'NEARProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.86'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +44.58%  '

# Row 36
$ws.Range("B36").Value = 'Dai'
$ws.Range("C36").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.04%  '

# Row 37
$ws.Range("E37").Value = '  +3.10%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '56.34'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +6.54%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0474'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +10.28%  '

# Row 40
$ws.Range("B40").Value = 'Fetch.AI'
$ws.Range("C40").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.60'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +51.71%  '

# Row 41
$ws.Range("B41").Value = 'PEPE'
$ws.Range("C41").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0₃0681'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +12.08%  '

# Row 42
$ws.Range("B42").Value = 'Stellar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.141'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +9.35%  '

# Row 43
$ws.Range("B43").Value = 'ThetaToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.83'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +7.75%  '

# Row 44
$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.34%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.37'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +11.73%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.321'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +18.67%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.13'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.82%  '

# Row 48
$ws.Range("E48").Value = '  +7.57%  '

# Row 49
$ws.Range("E49").Value = '  +8.53%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '142.52'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.80%  '

# Row 51
$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.81'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +8.18%  '

Write-Output "Applied cryptos update"
